# Amélioration de l'affichage des statistiques
# Met à jour le tableau "Contenu du stage" (langages utilisés en stage)
# avec les nouveaux effectifs et pourcentages calculés.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# --- Nombre d'étudiants par langage (colonne E) ---
$ws.Range("E16").Value = 13   # C#
$ws.Range("E17").Value = 28   # COBOL
$ws.Range("E19").Value = 3    # ASSEMBLEUR
$ws.Range("E20").Value = 1    # ANDROID

# --- Pourcentages correspondants (colonne G), saisis comme texte ---
# On force le format Texte avant la saisie pour éviter qu'Excel ne
# convertisse automatiquement la chaîne en valeur numérique pourcentage,
# puis on réinitialise la mise en forme de la cellule.

$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "28.89 %"
$ws.Range("G16").ClearFormats()

$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "62.22 %"
$ws.Range("G17").ClearFormats()

$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "6.67 %"
$ws.Range("G19").ClearFormats()

$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "2.22 %"
$ws.Range("G20").ClearFormats()

# --- Rafraîchit le graphique "Contenu du stage" pour qu'il reprenne
#     les nouvelles valeurs de la plage E16:E23 ---
$chartObj = $ws.ChartObjects().Item(2)
$chartObj.Chart.Refresh()
